# Fix full size and AF output for weather clusters
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Joint regime area" rows (36-40) that no longer belong in the
# 10-cluster weather assignment output (dimension shrinks from I40 to J35).
$ws.Rows("36:40").Delete()

# Add the new 9th cluster column (index 8) with its header value, copying the
# existing header formatting (bold/centered/bordered style) from column I.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = 8

# Updated cluster-assignment fractions per country (columns B:J, rows 2:35).
$ws.Range("D2").Value = 0
$ws.Range("F2").Value = 0.02010323281716925
$ws.Range("J2").Value = 0
$ws.Range("B3").Value = 0.2808564231738034
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0.007266629401900505
$ws.Range("F3").Value = 0.03069817984243413
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.08728010825439783
$ws.Range("J3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0.06976744186046512
$ws.Range("D4").Value = 0.03859250851305335
$ws.Range("E4").Value = 0
$ws.Range("G4").Value = 0.01875
$ws.Range("H4").Value = 0.138702460850112
$ws.Range("J4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.05254332029066506
$ws.Range("F5").Value = 0.04808475957620192
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0.07713125845737477
$ws.Range("J5").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.003912800447177195
$ws.Range("F6").Value = 0.02770986145069275
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("D7").Value = 0.1384790011350738
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("J7").Value = 0.01449275362318841
$ws.Range("B8").Value = 0.2122166246851379
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.04695360536612624
$ws.Range("F8").Value = 0.04373811464275997
$ws.Range("I8").Value = 0.08525033829499322
$ws.Range("J8").Value = 0
$ws.Range("D9").Value = 0.02497162315550511
$ws.Range("E9").Value = 0.02124091671324763
$ws.Range("F9").Value = 0
$ws.Range("J9").Value = 0.3864734299516917
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0.1162790697674418
$ws.Range("D10").Value = 0.1475595913734393
$ws.Range("E10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0.02684563758389261
$ws.Range("J10").Value = 0
$ws.Range("E11").Value = 0.01453325880380101
$ws.Range("F11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("D12").Value = 0.02043132803632236
$ws.Range("E12").Value = 0.0150922302962549
$ws.Range("F12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0.03542673107890499
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.00558971492453885
$ws.Range("F13").Value = 0.01358326541700625
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0.023906547133931
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("D16").Value = 0.1225879682179342
$ws.Range("E16").Value = 0.1704863051984343
$ws.Range("F16").Value = 0.03259983700081499
$ws.Range("G16").Value = 0.6281250000000005
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("F17").Value = 0.0271665308340125
$ws.Range("J17").Value = 0
$ws.Range("D18").Value = 0.01929625425652667
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("J18").Value = 0.001610305958132045
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0.01788708775852432
$ws.Range("F19").Value = 0.01059494702526488
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0.00111794298490777
$ws.Range("F20").Value = 0.1246943765281162
$ws.Range("J20").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0.04918949133594177
$ws.Range("F21").Value = 0.115186090736212
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("B23").Value = 0.1630982367758187
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0.04751257685858012
$ws.Range("F23").Value = 0.09236620483564174
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0.1204330175913398
$ws.Range("J23").Value = 0
$ws.Range("D24").Value = 0.03859250851305335
$ws.Range("E24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("E25").Value = 0.04415874790385683
$ws.Range("F25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0.07407407407407403
$ws.Range("E26").Value = 0.0005589714924538849
$ws.Range("F26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("D27").Value = 0.01929625425652667
$ws.Range("E27").Value = 0.02012297372833986
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0.003125
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0.008384572386808275
$ws.Range("F28").Value = 0.02662320021733225
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0.03186137506987144
$ws.Range("F29").Value = 0.001901657158380875
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("D30").Value = 0.001135073779795687
$ws.Range("E30").Value = 0.05980994969256552
$ws.Range("F30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("D31").Value = 0.04426787741203178
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("J31").Value = 0.1191626409017711
$ws.Range("B32").Value = 0
$ws.Range("C32").Value = 0.8139534883720939
$ws.Range("D32").Value = 0.09080590238365494
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0.3499999999999993
$ws.Range("H32").Value = 0.7472035794183463
$ws.Range("J32").Value = 0
$ws.Range("D33").Value = 0.2145289443813848
$ws.Range("E33").Value = 0.08719955282280603
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 0.08724832214765106
$ws.Range("J33").Value = 0
$ws.Range("D34").Value = 0
$ws.Range("F34").Value = 0.01168160825862538
$ws.Range("J34").Value = 0
$ws.Range("D35").Value = 0.01475595913734393
$ws.Range("E35").Value = 0.00111794298490777
$ws.Range("F35").Value = 0
$ws.Range("J35").Value = 0.07568438003220607
